$wb = $excel.ActiveWorkbook

# --- re_profiles sheet: header / label changes -----------------------------
$ws = $wb.Worksheets.Item("re_profiles")

# Order matters: new shared strings are appended in the order they are
# written, so "flo_fr" must be written before the "~TFM_DINS-AT: limtype=UP"
# label to match the target shared-string table ordering.
$ws.Range("I3").Value = "flo_fr"
$ws.Range("H2").Value = "~TFM_DINS-AT: limtype=UP"
$ws.Range("J3").Value = "pset_ci"

# J4 gets its own (non-shared) formula; J5:J68 and J69:J123 become two
# separate shared-formula groups, each simply mirroring column K.
$ws.Range("J4").Formula = "=K4"
$ws.Range("J5:J68").Formula = "=K5"
$ws.Range("J69:J123").Formula = "=K69"

# Widen columns J & K to fit the new content.
$ws.Columns("J:K").ColumnWidth = 10.166666666666666

# Make re_profiles the active sheet/tab, with H3 selected.
$ws.Activate()
[void]$ws.Range("H3").Select()
